$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:222 down to 104:223
$ws.Rows.Item(103).Insert()

# Populate the new row 103 with data
$ws.Range("A103").Value = 7
$ws.Range("B103").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C103").Value = "Ñuble"
$ws.Range("D103").Value = 45280
$ws.Range("E103").Value = 16
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100108
$ws.Range("H103").Value = "Tropicales y subtropicales"
$ws.Range("I103").Value = 100108002
$ws.Range("J103").Value = "Mango"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 100
$ws.Range("N103").Value = 10000
$ws.Range("O103").Value = 10000
$ws.Range("P103").Value = 10000
$ws.Range("Q103").Value = "$/bandeja 4 kilos"
$ws.Range("R103").Value = "Perú"
$ws.Range("S103").Value = 2500
$ws.Range("T103").Value = 4
